$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(608).Insert()

$ws.Cells.Item(608, 1).Value = "'2026/01/12"
$ws.Cells.Item(608, 1).Style = "Normal"
$ws.Cells.Item(608, 2).Value = "月"
$ws.Cells.Item(608, 3).Value = 13
$ws.Cells.Item(608, 4).Value = 21
